$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.102351
$ws.Range("H2").Value = 0.307053
$ws.Range("I2").Value = 0.2080046986044413
$ws.Range("J2").Value = 0.2080046986044413
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.005328
$ws.Range("N2").Value = 0.015984
$ws.Range("O2").Value = 0.001689940172269439
$ws.Range("P2").Value = 0.001689940172269439
$ws.Range("Q2").Value = 0.0005453261280000001
$ws.Range("R2").Value = 0.004907935152000001
$ws.Range("S2").Value = 0.0003515154961924424
$ws.Range("T2").Value = 0.0003515154961924424
$ws.Range("G3").Value = 0.102351
$ws.Range("H3").Value = 0.307053
$ws.Range("I3").Value = 0.2080046986044413
$ws.Range("J3").Value = 0.2080046986044413
$ws.Range("O3").Value = 0.9983100598277306
$ws.Range("P3").Value = 0.9983100598277306
$ws.Range("Q3").Value = 0.322144279663
$ws.Range("R3").Value = 2.899298516967
$ws.Range("S3").Value = 0.2076531831082489
$ws.Range("T3").Value = 0.2076531831082489
$ws.Range("I4").Value = 0.7162147240552154
$ws.Range("J4").Value = 0.7162147240552154
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.005328
$ws.Range("N4").Value = 0.015984
$ws.Range("O4").Value = 0.001689940172269439
$ws.Range("P4").Value = 0.001689940172269439
$ws.Range("Q4").Value = 0.001877700864
$ws.Range("R4").Value = 0.016899307776
$ws.Range("S4").Value = 0.00121036003415178
$ws.Range("T4").Value = 0.00121036003415178
$ws.Range("I5").Value = 0.7162147240552154
$ws.Range("J5").Value = 0.7162147240552154
$ws.Range("O5").Value = 0.9983100598277306
$ws.Range("P5").Value = 0.9983100598277306
$ws.Range("S5").Value = 0.7150043640210636
$ws.Range("T5").Value = 0.7150043640210636
$ws.Range("G6").Value = 0.03728866666666666
$ws.Range("I6").Value = 0.07578057734034331
$ws.Range("J6").Value = 0.0757805773403433
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.005328
$ws.Range("N6").Value = 0.015984
$ws.Range("O6").Value = 0.001689940172269439
$ws.Range("P6").Value = 0.001689940172269439
$ws.Range("Q6").Value = 0.000198674016
$ws.Range("R6").Value = 0.001788066144
$ws.Range("S6").Value = 0.0001280646419252174
$ws.Range("T6").Value = 0.0001280646419252173
$ws.Range("G7").Value = 0.03728866666666666
$ws.Range("I7").Value = 0.07578057734034331
$ws.Range("J7").Value = 0.0757805773403433
$ws.Range("O7").Value = 0.9983100598277306
$ws.Range("P7").Value = 0.9983100598277306
$ws.Range("S7").Value = 0.07565251269841809
$ws.Range("T7").Value = 0.07565251269841808
